# Add a new "TXT" output-format column to Sheet2, shifting the existing
# "Transpose" column one position to the right, and populate the new
# sample data row accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Previously M1 held "Transpose"; insert "TXT" in its place and push
# "Transpose" out to the new column N.
$ws.Range("N1").Value = "Transpose"
$ws.Range("M1").Value = "TXT"

# Fill in the new row-2 sample values.
$ws.Range("C2").Value = 41
$ws.Range("D2").Value = 3
$ws.Range("L2").Value = "No"
$ws.Range("N2").Value = "Yes"

# Match the selection left behind by the edit.
$ws.Range("N2").Select() | Out-Null
